$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("E2").Value = "Sitaram@55"
$ws.Range("D7").Select()
